$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Authentication codes for E2:E103, keyed by row number
$authCodes = @{
    2 = "980380"
    3 = "869215"
    4 = "962236"
    5 = "464023"
    6 = "936997"
    7 = "787965"
    8 = "994515"
    9 = "346857"
    10 = "730438"
    11 = "389118"
    12 = "729012"
    13 = "519737"
    14 = "385432"
    15 = "953153"
    16 = "628249"
    17 = "444923"
    18 = "176475"
    19 = "139507"
    20 = "738484"
    21 = "834355"
    22 = "864841"
    23 = "171803"
    24 = "181586"
    25 = "155354"
    26 = "213473"
    27 = "546779"
    28 = "768395"
    29 = "485369"
    30 = "585172"
    31 = "903143"
    32 = "173650"
    33 = "208394"
    34 = "768612"
    35 = "907900"
    36 = "704897"
    37 = "597581"
    38 = "854670"
    39 = "745577"
    40 = "612452"
    41 = "898091"
    42 = "126394"
    43 = "418642"
    44 = "541281"
    45 = "712603"
    46 = "910607"
    47 = "455458"
    48 = "253530"
    49 = "839338"
    50 = "109671"
    51 = "304480"
    52 = "222609"
    53 = "349320"
    54 = "538007"
    55 = "266802"
    56 = "889105"
    57 = "370184"
    58 = "358021"
    59 = "507495"
    60 = "396873"
    61 = "678386"
    62 = "604968"
    63 = "231300"
    64 = "375395"
    65 = "312840"
    66 = "149778"
    67 = "675561"
    68 = "985239"
    69 = "305000"
    70 = "491445"
    71 = "217039"
    72 = "962823"
    73 = "814247"
    74 = "662672"
    75 = "737056"
    76 = "106699"
    77 = "728342"
    78 = "989269"
    79 = "948420"
    80 = "183149"
    81 = "691711"
    82 = "497244"
    83 = "656516"
    84 = "931834"
    85 = "208546"
    86 = "803195"
    87 = "136982"
    88 = "405446"
    89 = "132462"
    90 = "411804"
    91 = "367603"
    92 = "258910"
    93 = "997294"
    94 = "647916"
    95 = "363254"
    96 = "247307"
    97 = "941276"
    98 = "763535"
    99 = "534988"
    100 = "785857"
    101 = "459648"
    102 = "695374"
    103 = "789268"
}

# Format the target range as Text first so the numeric-looking codes
# (which may contain leading zeros) are stored as text, matching the
# column's original text-cell type.
$codeRange = $ws.Range("E2:E103")
$codeRange.NumberFormat = "@"

foreach ($row in $authCodes.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $authCodes[$row]
}

